$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Step 1: paragraphs 1 and 2 lose their paragraph-mark rFonts/pPr ---
$p1 = $d.Paragraphs.Item(1)
$p2 = $d.Paragraphs.Item(2)
$rng12 = $d.Range($p1.Range.Start, $p2.Range.End)
$xml12 = "<w:p $wns><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>一、解决的问题</w:t></w:r></w:p>" + `
         "<w:p $wns><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>由于不同的机器有不同的操作系统，以及不同的库和组件，在将一个应用部署到多台机器上需要进行大量的环境配置操作。</w:t></w:r></w:p>"
[void]$rng12.InsertXML($xml12)

# --- Step 2: insert the new "virtual machine" section after the Docker paragraph ---
$p3 = $d.Paragraphs.Item(3)
$insPos = $p3.Range.End - 1
$insRng = $d.Range($insPos, $insPos)
$xmlNew = "<w:p $wns><w:pPr><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr></w:p>" + `
          "<w:p $wns><w:pPr><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>二、与虚拟机的比较</w:t></w:r></w:p>" + `
          "<w:p $wns><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>虚拟机也是一种虚拟化技术，它与</w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t xml:space=`"preserve`"> Docker </w:t></w:r><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>最大的区别在于它是通过模拟硬件，并在硬件上安装操作系统来实现。</w:t></w:r></w:p>"
[void]$insRng.InsertXML($xmlNew)
